# Update countries & provincias Spain
# Applies the diff: updates the "Datos actualizados" timestamp,
# updates COVID stats for several countries, and swaps the rank
# positions of Japon/Chile and Birmania/Congo/Barbados.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 06:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 614211
$ws.Range("C4").Value = 325
$ws.Range("E4").Value = 549327
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 26064

# --- Rows 27/28: Japon moves above Chile (Japon updated, Chile unchanged but shifted down) ---
$ws.Range("A27").Value = "Japon"
$ws.Range("B27").Value = 8100
$ws.Range("C27").Value = 215
$ws.Range("D27").Value = 853
$ws.Range("E27").Value = 7101
$ws.Range("F27").Value = 152
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 146

$ws.Range("A28").Value = "Chile"
$ws.Range("B28").Value = 7917
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 2646
$ws.Range("E28").Value = 5179
$ws.Range("F28").Value = 387
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 92

# --- Row 34: Australia ---
$ws.Range("B34").Value = 6440
$ws.Range("C34").Value = 40
$ws.Range("E34").Value = 2779
$ws.Range("F34").Value = 78
$ws.Range("G34").Value = 2
$ws.Range("H34").Value = 63

# --- Row 36: Pakistan ---
$ws.Range("B36").Value = 5983
$ws.Range("C36").Value = 146
$ws.Range("E36").Value = 4498
$ws.Range("G36").Value = 11
$ws.Range("H36").Value = 107

# --- Row 69: Kazajistan ---
$ws.Range("B69").Value = 1267
$ws.Range("C69").Value = 35
$ws.Range("E69").Value = 1050
$ws.Range("F69").Value = 20

# --- Row 120: Venezuela ---
$ws.Range("B120").Value = 197
$ws.Range("C120").Value = 8
$ws.Range("E120").Value = 77

# --- Rows 141/142/143: Birmania moves above Congo/Barbados ---
$ws.Range("A141").Value = "Birmania"
$ws.Range("B141").Value = 74
$ws.Range("C141").Value = 11
$ws.Range("D141").Value = 2
$ws.Range("E141").Value = 68
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 4

$ws.Range("A142").Value = "Congo"
$ws.Range("B142").Value = 74
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 10
$ws.Range("E142").Value = 59
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 5

$ws.Range("A143").Value = "Barbados"
$ws.Range("B143").Value = 73
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 15
$ws.Range("E143").Value = 53
$ws.Range("F143").Value = 4
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 5
